$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number + date range) ---
$ws.Range("A8").Value = "Volume 31   Number  52"
$ws.Range("C9").Value = "Report Covering the Week  12/23/2024  Through  12/29/2024"

# --- Convert numeric cells to text placeholders ("0" / "***.*") using donor cells with matching style ---
$ws.Range("C14").Copy($ws.Range("C18"))
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("C14").Copy($ws.Range("G22"))
$ws.Range("E14").Copy($ws.Range("H22"))
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))

# --- Convert text placeholder cells to numeric (copy donor style, then set value) ---
$ws.Range("C16").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2

# --- Plain numeric value updates ---
$ws.Range("F16").Value = 3
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = -62.5
$ws.Range("I16").Value = 89
$ws.Range("J16").Value = 124
$ws.Range("K16").Value = -28.225806451612
$ws.Range("L16").Value = -36.428571428571
$ws.Range("M16").Value = -31.538461538461
$ws.Range("N16").Value = -85.067114093959
$ws.Range("C17").Value = 2
$ws.Range("D17").Value = 4
$ws.Range("F17").Value = 8
$ws.Range("H17").Value = 14.285714285714
$ws.Range("I17").Value = 129
$ws.Range("J17").Value = 115
$ws.Range("K17").Value = 12.173913043478
$ws.Range("L17").Value = -17.307692307692
$ws.Range("M17").Value = 76.712328767123
$ws.Range("N17").Value = -51.503759398496
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("F18").Value = 6
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = -50
$ws.Range("J18").Value = 174
$ws.Range("K18").Value = -17.816091954023
$ws.Range("L18").Value = -46.240601503759
$ws.Range("M18").Value = -40.663900414937
$ws.Range("N18").Value = -87.266251113089
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 10
$ws.Range("E19").Value = -30
$ws.Range("F19").Value = 40
$ws.Range("G19").Value = 50
$ws.Range("H19").Value = -20
$ws.Range("I19").Value = 666
$ws.Range("J19").Value = 681
$ws.Range("K19").Value = -2.202643171806
$ws.Range("L19").Value = 9.719934102141
$ws.Range("M19").Value = 114.147909967846
$ws.Range("N19").Value = 80.978260869565
$ws.Range("C20").Value = 1
$ws.Range("G20").Value = 6
$ws.Range("H20").Value = 33.333333333333
$ws.Range("I20").Value = 112
$ws.Range("K20").Value = -26.797385620915
$ws.Range("L20").Value = -41.666666666666
$ws.Range("M20").Value = -22.758620689655
$ws.Range("N20").Value = -87.732749178532
$ws.Range("C21").Value = 11
$ws.Range("D21").Value = 19
$ws.Range("E21").Value = -42.105263157894
$ws.Range("F21").Value = 68
$ws.Range("G21").Value = 83
$ws.Range("H21").Value = -18.072289156626
$ws.Range("I21").Value = 1152
$ws.Range("J21").Value = 1253
$ws.Range("K21").Value = -8.060654429369
$ws.Range("L21").Value = -16.279069767441
$ws.Range("M21").Value = 27.857935627081
$ws.Range("N21").Value = -64.89945155393
$ws.Range("M22").Value = -5.882352941176
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = 28.571428571428
$ws.Range("F24").Value = 70
$ws.Range("G24").Value = 65
$ws.Range("H24").Value = 7.692307692307
$ws.Range("I24").Value = 985
$ws.Range("J24").Value = 897
$ws.Range("K24").Value = 9.810479375696
$ws.Range("L24").Value = -1.696606786427
$ws.Range("M24").Value = 61.740558292282
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = 125
$ws.Range("F25").Value = 32
$ws.Range("G25").Value = 35
$ws.Range("H25").Value = -8.571428571428
$ws.Range("I25").Value = 611
$ws.Range("J25").Value = 492
$ws.Range("K25").Value = 24.186991869918
$ws.Range("L25").Value = 3.91156462585
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = 200
$ws.Range("F26").Value = 17
$ws.Range("G26").Value = 13
$ws.Range("H26").Value = 30.76923076923
$ws.Range("I26").Value = 277
$ws.Range("J26").Value = 238
$ws.Range("K26").Value = 16.386554621848
$ws.Range("L26").Value = -0.716845878136
$ws.Range("M26").Value = 34.466019417475
$ws.Range("E28").Value = 100
$ws.Range("F28").Value = 3
$ws.Range("G28").Value = 3
$ws.Range("H28").Value = 0
$ws.Range("I28").Value = 42
$ws.Range("J28").Value = 45
$ws.Range("K28").Value = -6.666666666666
$ws.Range("L28").Value = 44.827586206896
